$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 1: drop the redundant extra "value" header cells (C1:F1), leaving A1:B1
$ws.Range("C1:F1").ClearContents()

# L_curve parameter value changes from 0 to 1
$ws.Range("B9").Value = 1

# Remove the obsolete "Deletion" parameter row entirely (rows below shift up)
$ws.Rows.Item(17).Delete()
